$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------------
# Row 1 (original): "1 <tab> 0.00326 <tab> 0.00326 <tab> 0.00326 <tab> 0.00000
#   <tab> 0.00326 <tab> 0.00326 <tab> 0.00326 <tab> 0.00326 <tab> 100.0"
# becomes 13 single-value rows.
# ---------------------------------------------------------------------------
$row1Values = @("0M","0M","0M","1682","0.00003","0.25696","0.04413","0.03404","0.02351","0.13598","0.19359","104.09665","100.0")

# Set the first (existing) row's cell to the first new value.
$t.Rows.Item(1).Cells.Item(1).Range.Text = $row1Values[0]

# Insert the remaining 12 values as new rows right after row 1, by repeatedly
# inserting immediately before the row that currently follows row 1 (which,
# at the start of this loop, is the original row 2). Walking the replacement
# values in reverse and always inserting before that same anchor row yields
# the values in ascending (correct) order once the loop completes.
$anchorAfterRow1 = $t.Rows.Item(2)
for ($i = $row1Values.Length - 1; $i -ge 1; $i--) {
    $newRow = $t.Rows.Add($anchorAfterRow1)
    $newRow.Cells.Item(1).Range.Text = $row1Values[$i]
}

# ---------------------------------------------------------------------------
# Row 2 (original, now at table position 14): "153 <tab> 0.00027 <tab> ...
#   <tab> 100.0" becomes a single row with text "0".
# ---------------------------------------------------------------------------
$t.Rows.Item(14).Cells.Item(1).Range.Text = "0"

# ---------------------------------------------------------------------------
# Row 3 (original, now at table position 15): "142 <tab> 0.00088 <tab> ...
#   <tab> 100.0" becomes 16 single-value rows.
# ---------------------------------------------------------------------------
$row3Values = @("0.00000","0.00000","0.00000","0.00000","0.00000","0.00000","0.00000","0.00000","0.0","0","0.00000","0.00000","0.00000","0.00000","0.00000","0.0")

$t.Rows.Item(15).Cells.Item(1).Range.Text = $row3Values[0]

$anchorAfterRow3 = $t.Rows.Item(16)
for ($i = $row3Values.Length - 1; $i -ge 1; $i--) {
    $newRow = $t.Rows.Add($anchorAfterRow3)
    $newRow.Cells.Item(1).Range.Text = $row3Values[$i]
}

# ---------------------------------------------------------------------------
# Rows that originally held 1841, 0.00000, 0.50118, 0.02832, 0.01281,
# 0.08405, 0.16850, 0.17257, 52.17620, 50.1 are left untouched. After the
# insertions above, the table now has 13 + 1 + 16 = 30 rows standing in for
# the original first three rows, so those unchanged rows now sit at
# positions 31-40.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Original row 14 (ten tab-separated zeros), now at position 41, becomes a
# single run with text "-66.86".
# ---------------------------------------------------------------------------
$t.Rows.Item(41).Cells.Item(1).Range.Text = "-66.86"

# ---------------------------------------------------------------------------
# Original row 15 (empty run), now at position 42, gains the text "104.1".
# ---------------------------------------------------------------------------
$t.Rows.Item(42).Cells.Item(1).Range.Text = "104.1"

# ---------------------------------------------------------------------------
# Original row 16 (ten tab-separated zeros), now at position 43, becomes a
# single run with text "62".
# ---------------------------------------------------------------------------
$t.Rows.Item(43).Cells.Item(1).Range.Text = "62"
